# Adds Output [Data] / Data Format / Data Selector Format columns (with a
# second data row) to the existing mzML and mzMLite annotation tables, and
# adds a new "psmstats" sheet + annotation table that carries the mzlite ->
# psm step of the pipeline.

$wb = $excel.ActiveWorkbook

$wiff1   = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/assays/dilutionSeries_Chlamy/dataset/20170519 TM FScon3501/20170519 TM FScon3501.wiff"
$wiff2   = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/assays/dilutionSeries_Chlamy/dataset/20170519 TM FScon3501/20170519 TM FScon3503.wiff"
$mzml1   = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/mzml/20170519_TM_FScon3501.mzML"
$mzml2   = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/mzml/20170519_TM_FScon3503.mzML"
$mzlite1 = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/mzlite/20170519 TM FScon3501.mzlite"
$mzlite2 = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/mzlite/20170519 TM FScon3503.mzlite"
$psm1    = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/psm/20170519 TM FScon3501.psm"
$psm2    = "/home/paulinehans/Dokumente/QualIQon_all/QualIQon_v1.0/runs/dilution_series/psm/20170519 TM FScon3503.psm"

function Expand-AnnotationTable {
    param($ws, $inA2, $outD2, $inA3, $outD3)

    $lo = $ws.ListObjects.Item(1)

    # 3 new columns: Output [Data] / Data Format  / Data Selector Format
    $lo.ListColumns.Add() | Out-Null
    $lo.ListColumns.Add() | Out-Null
    $lo.ListColumns.Add() | Out-Null
    # 1 new row so the table holds two data rows
    $lo.ListRows.Add() | Out-Null

    $ws.Range("D1").Value = "Output [Data]"
    $ws.Range("E1").Value = "Data Format "
    $ws.Range("F1").Value = "Data Selector Format "

    $ws.Range("A2").Value = $inA2
    $ws.Range("D2").Value = $outD2
    # propagate the blank Data Format / Data Selector Format cells into the
    # newly added columns so row 2 stays fully populated
    $ws.Range("B2:C2").Copy($ws.Range("E2:F2"))

    $ws.Range("A3").Value = $inA3
    $ws.Range("D3").Value = $outD3
    $ws.Range("B2:C2").Copy($ws.Range("B3:C3"))
    $ws.Range("B2:C2").Copy($ws.Range("E3:F3"))
}

# --- mzML sheet: wiff (input) -> mzML (output) ---------------------------
$wsMzml = $wb.Worksheets.Item("mzML ")
Expand-AnnotationTable $wsMzml $wiff1 $mzml1 $wiff2 $mzml2

# --- mzMLite sheet: mzML (input) -> mzlite (output) -----------------------
$wsMzmlite = $wb.Worksheets.Item("mzMLite")
Expand-AnnotationTable $wsMzmlite $mzml1 $mzlite1 $mzml2 $mzlite2

# --- new psmstats sheet: mzlite (input) -> psm (output) -------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "psmstats"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# $newSheet becomes stale once Move() re-shuffles sheet positions, so grab a
# fresh reference by name before touching it again.
$wsPsm = $wb.Worksheets.Item("psmstats")

$wsPsm.Range("A1").Value = "Input [Data]"
$wsPsm.Range("B1").Value = "Data Format"
$wsPsm.Range("C1").Value = "Data Selector Format"

$loPsm = $wsPsm.ListObjects.Add(1, $wsPsm.Range("A1:C1"), $null, 1)
$loPsm.Name = "annotationTable2"

Expand-AnnotationTable $wsPsm $mzlite1 $psm1 $mzlite2 $psm2
